$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 43.841169
$ws.Range("H2").Value = 131.523507
$ws.Range("I2").Value = 0.6105408572336042
$ws.Range("J2").Value = 0.6105408572336042
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.972196
$ws.Range("N2").Value = 2.916588
$ws.Range("O2").Value = 0.7027023771175303
$ws.Range("P2").Value = 0.7027023771175303
$ws.Range("Q2").Value = 42.622209137124
$ws.Range("R2").Value = 383.599882234116
$ws.Range("S2").Value = 0.4290285117054283
$ws.Range("T2").Value = 0.4290285117054283

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 43.841169
$ws.Range("H3").Value = 131.523507
$ws.Range("I3").Value = 0.6105408572336042
$ws.Range("J3").Value = 0.6105408572336042
$ws.Range("O3").Value = 0.1592492623233027
$ws.Range("P3").Value = 0.1592492623233027
$ws.Range("Q3").Value = 9.659217877586999
$ws.Range("R3").Value = 86.93296089828299
$ws.Range("S3").Value = 0.09722818113268834
$ws.Range("T3").Value = 0.09722818113268834

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 43.841169
$ws.Range("H4").Value = 131.523507
$ws.Range("I4").Value = 0.6105408572336042
$ws.Range("J4").Value = 0.6105408572336042
$ws.Range("M4").Value = 0.1909913333333333
$ws.Range("N4").Value = 0.572974
$ws.Range("O4").Value = 0.138048360559167
$ws.Range("P4").Value = 0.138048360559167
$ws.Range("Q4").Value = 8.373283322202
$ws.Range("R4").Value = 75.35954989981799
$ws.Range("S4").Value = 0.0842841643954875
$ws.Range("T4").Value = 0.0842841643954875

# Row 5
$ws.Range("I5").Value = 0.0635739353967235
$ws.Range("J5").Value = 0.06357393539672351
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.972196
$ws.Range("N5").Value = 2.916588
$ws.Range("O5").Value = 0.7027023771175303
$ws.Range("P5").Value = 0.7027023771175303
$ws.Range("Q5").Value = 4.438133071759999
$ws.Range("R5").Value = 39.94319764584
$ws.Range("S5").Value = 0.0446735555259939
$ws.Range("T5").Value = 0.04467355552599391

# Row 6
$ws.Range("I6").Value = 0.0635739353967235
$ws.Range("J6").Value = 0.06357393539672351
$ws.Range("O6").Value = 0.1592492623233027
$ws.Range("P6").Value = 0.1592492623233027
$ws.Range("S6").Value = 0.01012410231491752
$ws.Range("T6").Value = 0.01012410231491752

# Row 7
$ws.Range("I7").Value = 0.0635739353967235
$ws.Range("J7").Value = 0.06357393539672351
$ws.Range("M7").Value = 0.1909913333333333
$ws.Range("N7").Value = 0.572974
$ws.Range("O7").Value = 0.138048360559167
$ws.Range("P7").Value = 0.138048360559167
$ws.Range("Q7").Value = 0.8718868961466666
$ws.Range("R7").Value = 7.84698206532
$ws.Range("S7").Value = 0.008776277555812076
$ws.Range("T7").Value = 0.008776277555812078

# Row 8
$ws.Range("G8").Value = 23.400872
$ws.Range("H8").Value = 70.202616
$ws.Range("I8").Value = 0.3258852073696723
$ws.Range("J8").Value = 0.3258852073696723
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.972196
$ws.Range("N8").Value = 2.916588
$ws.Range("O8").Value = 0.7027023771175303
$ws.Range("P8").Value = 0.7027023771175303
$ws.Range("Q8").Value = 22.750234154912
$ws.Range("R8").Value = 204.752107394208
$ws.Range("S8").Value = 0.229000309886108
$ws.Range("T8").Value = 0.229000309886108

# Row 9
$ws.Range("G9").Value = 23.400872
$ws.Range("H9").Value = 70.202616
$ws.Range("I9").Value = 0.3258852073696723
$ws.Range("J9").Value = 0.3258852073696723
$ws.Range("O9").Value = 0.1592492623233027
$ws.Range("P9").Value = 0.1592492623233027
$ws.Range("Q9").Value = 5.155750321656
$ws.Range("R9").Value = 46.401752894904
$ws.Range("S9").Value = 0.05189697887569685
$ws.Range("T9").Value = 0.05189697887569685

# Row 10
$ws.Range("G10").Value = 23.400872
$ws.Range("H10").Value = 70.202616
$ws.Range("I10").Value = 0.3258852073696723
$ws.Range("J10").Value = 0.3258852073696723
$ws.Range("M10").Value = 0.1909913333333333
$ws.Range("N10").Value = 0.572974
$ws.Range("O10").Value = 0.138048360559167
$ws.Range("P10").Value = 0.138048360559167
$ws.Range("Q10").Value = 4.469363744442667
$ws.Range("R10").Value = 40.224273699984
$ws.Range("S10").Value = 0.04498791860786743
$ws.Range("T10").Value = 0.04498791860786743

